$wb = $excel.ActiveWorkbook

# --- Filepath sheet: insert two new rows (raw generation/exchange directory) ---
$ws = $wb.Worksheets.Item("Filepath")
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "raw generation directory"
$ws.Range("A5").Value = "raw exchange directory"

# Filepath sheet is no longer the active one; reset its selection to A1
[void]$ws.Range("A1").Select()

# --- Parameter sheet becomes the active / selected sheet, with A1 selected ---
$wsParam = $wb.Worksheets.Item("Parameter")
$wsParam.Activate()
[void]$wsParam.Range("A1").Select()
